# Library management system added
# Adds three new dated entries (rows 34, 36, 38) to the learning log,
# mirroring the style/format of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entries: date (serial) in column A, notes text in column C.
$entries = @(
    @{ Row = 34; Date = 45497; Text = "depencency injection method=>constructor,property,function,webapi completed with database connection ,migration ,code first,db first approach." },
    @{ Row = 36; Date = 45498; Text = "web api using dbfirst approach ,data table created for the project,lms api project category post and books get,post api" },
    @{ Row = 38; Date = 45499; Text = "LibraryMgmtSystem project addbook,user,category,borrow record =>post,get api controller completed" }
)

# Use an existing date cell as the formatting source so the new date cells
# reuse the same (already-present) date style instead of creating a new one.
$dateFormatSource = $ws.Range("A32")

foreach ($entry in $entries) {
    $dateCell = $ws.Cells.Item($entry.Row, 1)
    $dateCell.Value = $entry.Date

    $dateFormatSource.Copy() | Out-Null
    $dateCell.PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($entry.Row, 3).Value = $entry.Text
}

# Match the final selection recorded in the workbook.
$ws.Range("C38").Select() | Out-Null
